$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2
$ws.Range("C2").Value = 0.5155555555555555
$ws.Range("J2").Value = 0.008888888888888889
$ws.Range("P2").Value = 0.1511111111111111
$ws.Range("S2").Value = 0.1244444444444444
$ws.Range("B3").Value = 0.007936507936507936
$ws.Range("C3").Value = 0.07142857142857142
$ws.Range("J3").Value = 0.02380952380952381
$ws.Range("P3").Value = 0.6904761904761905
$ws.Range("S3").Value = 0.2063492063492063
$ws.Range("J4").Value = 0.02857142857142857
$ws.Range("P4").Value = 0.6857142857142857
$ws.Range("S4").Value = 0.2857142857142857
$ws.Range("B6").Value = 0.02985074626865672
$ws.Range("D6").Value = 0.01119402985074627
$ws.Range("E6").Value = 0.003731343283582089
$ws.Range("F6").Value = 0.05970149253731343
$ws.Range("J6").Value = 0.25
$ws.Range("O6").Value = 0.02238805970149254
$ws.Range("Q6").Value = 0.1641791044776119
$ws.Range("R6").Value = 0.1082089552238806
$ws.Range("S6").Value = 0.3507462686567164
$ws.Range("B7").Value = 0.1031390134529148
$ws.Range("D7").Value = 0.03587443946188341
$ws.Range("F7").Value = 0.06278026905829596
$ws.Range("J7").Value = 0.1076233183856502
$ws.Range("O7").Value = 0.03139013452914798
$ws.Range("Q7").Value = 0.1569506726457399
$ws.Range("R7").Value = 0.09417040358744394
$ws.Range("S7").Value = 0.4080717488789238
$ws.Range("B8").Value = 0.06903765690376569
$ws.Range("D8").Value = 0.0104602510460251
$ws.Range("E8").Value = 0.002092050209205021
$ws.Range("F8").Value = 0.08577405857740586
$ws.Range("J8").Value = 0.09205020920502092
$ws.Range("O8").Value = 0.01882845188284519
$ws.Range("Q8").Value = 0.198744769874477
$ws.Range("R8").Value = 0.100418410041841
$ws.Range("S8").Value = 0.4225941422594142
$ws.Range("B9").Value = 0.06944444444444445
$ws.Range("D9").Value = 0.009259259259259259
$ws.Range("E9").Value = 0.004629629629629629
$ws.Range("F9").Value = 0.02777777777777778
$ws.Range("J9").Value = 0.07407407407407407
$ws.Range("O9").Value = 0.004629629629629629
$ws.Range("Q9").Value = 0.1759259259259259
$ws.Range("R9").Value = 0.1296296296296296
$ws.Range("S9").Value = 0.5046296296296297
$ws.Range("B10").Value = 0.08376963350785341
$ws.Range("D10").Value = 0.01657940663176265
$ws.Range("E10").Value = 0.0008726003490401396
$ws.Range("F10").Value = 0.0968586387434555
$ws.Range("J10").Value = 0.1029668411867365
$ws.Range("O10").Value = 0.02530541012216405
$ws.Range("Q10").Value = 0.1893542757417103
$ws.Range("R10").Value = 0.09162303664921466
$ws.Range("S10").Value = 0.3926701570680629
$ws.Range("G11").Value = 0.1496815286624204
$ws.Range("J11").Value = 0.06687898089171974
$ws.Range("K11").Value = 0.2006369426751592
$ws.Range("L11").Value = 0.5732484076433121
$ws.Range("S11").Value = 0.009554140127388535
$ws.Range("G12").Value = 0.7807486631016043
$ws.Range("J12").Value = 0.1711229946524064
$ws.Range("K12").Value = 0.0053475935828877
$ws.Range("L12").Value = 0.03208556149732621
$ws.Range("S12").Value = 0.0106951871657754
$ws.Range("G13").Value = 0.74
$ws.Range("J13").Value = 0.22
$ws.Range("S13").Value = 0.04
$ws.Range("F15").Value = 0.01626016260162602
$ws.Range("H15").Value = 0.1666666666666667
$ws.Range("I15").Value = 0.09349593495934959
$ws.Range("J15").Value = 0.3130081300813008
$ws.Range("K15").Value = 0.07317073170731707
$ws.Range("M15").Value = 0.008130081300813009
$ws.Range("N15").Value = 0.008130081300813009
$ws.Range("O15").Value = 0.05284552845528456
$ws.Range("S15").Value = 0.2682926829268293
$ws.Range("F16").Value = 0.02127659574468085
$ws.Range("H16").Value = 0.1914893617021277
$ws.Range("I16").Value = 0.05673758865248227
$ws.Range("J16").Value = 0.3829787234042553
$ws.Range("K16").Value = 0.07801418439716312
$ws.Range("M16").Value = 0.01418439716312057
$ws.Range("O16").Value = 0.09929078014184398
$ws.Range("S16").Value = 0.1560283687943262
$ws.Range("F17").Value = 0.02083333333333333
$ws.Range("H17").Value = 0.1875
$ws.Range("I17").Value = 0.1064814814814815
$ws.Range("J17").Value = 0.3703703703703703
$ws.Range("K17").Value = 0.09490740740740741
$ws.Range("M17").Value = 0.02546296296296296
$ws.Range("N17").Value = 0.002314814814814815
$ws.Range("O17").Value = 0.06712962962962964
$ws.Range("S17").Value = 0.125
$ws.Range("F18").Value = 0.02542372881355932
$ws.Range("H18").Value = 0.1991525423728814
$ws.Range("I18").Value = 0.1355932203389831
$ws.Range("J18").Value = 0.326271186440678
$ws.Range("K18").Value = 0.1186440677966102
$ws.Range("M18").Value = 0.01694915254237288
$ws.Range("O18").Value = 0.0635593220338983
$ws.Range("S18").Value = 0.1144067796610169
$ws.Range("F19").Value = 0.01920122887864823
$ws.Range("H19").Value = 0.2188940092165899
$ws.Range("I19").Value = 0.08141321044546851
$ws.Range("J19").Value = 0.3494623655913979
$ws.Range("K19").Value = 0.1205837173579109
$ws.Range("M19").Value = 0.02457757296466974
$ws.Range("N19").Value = 0.0007680491551459293
$ws.Range("O19").Value = 0.07450076804915515
$ws.Range("S19").Value = 0.1105990783410138
